$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Indicator text lost the period after "9.2.1"
$ws.Range("B4").Value = "9.2.1 Manufacturing value added as a proportion of GDP and per capita"

# Organization website updated to the gov.kg domain
$ws.Range("B10").Value = "www.stat.gov.kg"

# Selection moved from B2 to B9
$null = $ws.Range("B9").Select()
